$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1152
$ws.Range("B4").Value = 517
$ws.Range("B7").Value = 194
$ws.Range("B12").Value = 70
